# Replace every occurrence of $oldText with $newText while preserving the
# character formatting of the text being replaced. A naive
# Find.Execute(..., Replace:=2) / Range.Text= assignment in this runtime
# tends to pick up the formatting of a neighboring run and merge the
# replacement into it, which can silently drop character formatting
# (e.g. bold) from the replaced span. Inserting the new text at an
# interior point of the matched range (not exactly at its Start/End
# boundary) avoids that ambiguity, because an interior point unambiguously
# belongs to the run being replaced.
function Replace-All($doc, $oldText, $newText) {
    $count = 0
    while ($true) {
        $rng = $doc.Content
        $found = $rng.Find.Execute($oldText)
        if (-not $found) {
            break
        }

        $s = $rng.Start
        $e = $rng.End
        $oldLen = $e - $s
        $mid = $s + [int][Math]::Floor($oldLen / 2)

        # Insert the replacement text right after the interior midpoint;
        # it inherits the formatting of the character at $mid, which is
        # guaranteed to be inside the original match.
        $insPoint = $doc.Range($mid, $mid)
        $insPoint.InsertAfter($newText)
        $newLen = $newText.Length

        # Delete the trailing remainder of the old text (after the
        # inserted block), then the leading remainder (before it).
        $afterPart = $doc.Range($mid + $newLen, $e + $newLen)
        if ($afterPart.Start -lt $afterPart.End) {
            $afterPart.Delete()
        }
        $beforePart = $doc.Range($s, $mid)
        if ($beforePart.Start -lt $beforePart.End) {
            $beforePart.Delete()
        }

        $count = $count + 1
    }
    return $count
}

$d = $word.ActiveDocument

# Party name: WEBB WINTERS TRADING -> RAMOS AND ROY PLC
# (handle the ", " suffixed occurrence first so the plain search string
# below does not also match inside it)
Replace-All $d "WEBB WINTERS TRADING ," "RAMOS AND ROY PLC ,"
Replace-All $d "WEBB WINTERS TRADING" "RAMOS AND ROY PLC"

# Party name: Lawrence Morse Traders -> Anthony And Serrano Llc
Replace-All $d "Lawrence Morse Traders " "Anthony And Serrano Llc "

# Signatory name: MERRILL GILLESPIE -> KYLEE MORALES
# (trailing-space variant first, for the same reason as above)
Replace-All $d "MERRILL GILLESPIE " "KYLEE MORALES "
Replace-All $d "MERRILL GILLESPIE" "KYLEE MORALES"

# Dates: 26 February 2025/2027 -> 01 March 2025/2027
Replace-All $d "26 February 2025" "01 March 2025"
Replace-All $d "26 February 2027 " "01 March 2027 "

# Name: LOIS ESTES -> LUKE MARSH
Replace-All $d "LOIS ESTES" "LUKE MARSH"

# Position: ELIGENDI ARCHITECTO -> INCIDUNT IURE DOLOR
Replace-All $d "ELIGENDI ARCHITECTO " "INCIDUNT IURE DOLOR"
